# Update the lattice-multiplication exercise table: each of the 15 cells
# holds a single run with 5 lines (problem, top factors, "----", two
# left-side digits) separated by <w:br/>. Word exposes line breaks inside
# Range.Text as vertical-tab (chr 11) characters, so rewriting Cell.Range.Text
# with vt-joined strings replaces the run's text while preserving the run's
# formatting (rPr sz=32) and the <w:t>/<w:br/> structure.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$c = $t.Cell(1,1)
$c.Range.Text = "59 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "5|    |" + $vt + "9|    |"

$c = $t.Cell(1,2)
$c.Range.Text = "30 x 86" + $vt + "  8    6" + $vt + "  ----" + $vt + "3|    |" + $vt + "0|    |"

$c = $t.Cell(1,3)
$c.Range.Text = "75 x 57" + $vt + "  5    7" + $vt + "  ----" + $vt + "7|    |" + $vt + "5|    |"

$c = $t.Cell(2,1)
$c.Range.Text = "65 x 56" + $vt + "  5    6" + $vt + "  ----" + $vt + "6|    |" + $vt + "5|    |"

$c = $t.Cell(2,2)
$c.Range.Text = "17 x 70" + $vt + "  7    0" + $vt + "  ----" + $vt + "1|    |" + $vt + "7|    |"

$c = $t.Cell(2,3)
$c.Range.Text = "17 x 30" + $vt + "  3    0" + $vt + "  ----" + $vt + "1|    |" + $vt + "7|    |"

$c = $t.Cell(3,1)
$c.Range.Text = "56 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "5|    |" + $vt + "6|    |"

$c = $t.Cell(3,2)
$c.Range.Text = "45 x 53" + $vt + "  5    3" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"

$c = $t.Cell(3,3)
$c.Range.Text = "42 x 42" + $vt + "  4    2" + $vt + "  ----" + $vt + "4|    |" + $vt + "2|    |"

$c = $t.Cell(4,1)
$c.Range.Text = "46 x 23" + $vt + "  2    3" + $vt + "  ----" + $vt + "4|    |" + $vt + "6|    |"

$c = $t.Cell(4,2)
$c.Range.Text = "47 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "4|    |" + $vt + "7|    |"

$c = $t.Cell(4,3)
$c.Range.Text = "49 x 70" + $vt + "  7    0" + $vt + "  ----" + $vt + "4|    |" + $vt + "9|    |"

$c = $t.Cell(5,1)
$c.Range.Text = "18 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "1|    |" + $vt + "8|    |"

$c = $t.Cell(5,2)
$c.Range.Text = "53 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "5|    |" + $vt + "3|    |"

$c = $t.Cell(5,3)
$c.Range.Text = "83 x 10" + $vt + "  1    0" + $vt + "  ----" + $vt + "8|    |" + $vt + "3|    |"
